# Fall 22 week 14 complete - append 25 new matchup rows (266-290) to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, 17, 5, 3),
    @(4, 8, 6, 12),
    @(4, 4, 3, 16),
    @(1, 8, 3, 12),
    @(5, 7, 4, 13),
    @(5, 8, 4, 12),
    @(3, 8, 4, 12),
    @(7, 15, 4, 5),
    @(6, 14, 5, 6),
    @(3, 15, 4, 5),
    @(6, 13, 5, 7),
    @(3, 15, 5, 5),
    @(3, 14, 2, 6),
    @(2, 7, 4, 13),
    @(3, 12, 4, 8),
    @(8, 12, 7, 8),
    @(5, 5, 2, 15),
    @(3, 14, 4, 6),
    @(4, 15, 3, 5),
    @(4, 7, 6, 13),
    @(5, 14, 6, 6),
    @(8, 15, 9, 5),
    @(4, 4, 5, 16),
    @(1, 12, 2, 8),
    @(5, 14, 3, 6)
)

$startRow = 266
$row = $startRow
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Formula = "=B$row+D$row"
    $row = $row + 1
}

$lastRow = $row - 1

# Update the active selection / view to match the new bottom of data
[void]$ws.Range("A$($lastRow + 1)").Select()
